$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.977634906768799
$ws.Range("B1").Value = 5.096848011016846
$ws.Range("C1").Value = 4.194746971130371
$ws.Range("D1").Value = 4.952110290527344
$ws.Range("E1").Value = 5.269463062286377
